$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Cell updates derived from the cryptos list refresh (prices + 1h volume %)
$ws.Range("D2").Value = "42.992.91"
$ws.Range("E2").Value = "  +4.24%  "

$ws.Range("D3").Value = "2.232.74"
$ws.Range("E3").Value = "  +4.12%  "

$ws.Range("E4").Value = "  -0.10%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "251.92"
$ws.Range("E5").Value = "  +6.53%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.608"
$ws.Range("E6").Value = "  +1.18%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "75.19"
$ws.Range("E7").Value = "  +7.77%  "

$ws.Range("E8").Value = "  -0.14%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.596"
$ws.Range("E9").Value = "  +4.57%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "41.04"
$ws.Range("E10").Value = "  +4.97%  "

$ws.Range("E11").Value = "  +3.14%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "6.88"
$ws.Range("E12").Value = "  +4.15%  "

$ws.Range("E13").Value = "  +2.28%  "

$ws.Range("D14").Value = "2.570.64"
$ws.Range("E14").Value = "  +4.28%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "14.44"
$ws.Range("E15").Value = "  +0.97%  "

$ws.Range("D16").Value = "2.230.44"
$ws.Range("E16").Value = "  +4.46%  "

$ws.Range("E17").Value = "  +1.46%  "

$ws.Range("D18").Value = "42.881.75"
$ws.Range("E18").Value = "  +4.46%  "

$ws.Range("E19").Value = "  +3.78%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "71.16"
$ws.Range("E20").Value = "  +3.37%  "

$ws.Range("E21").Value = "  +4.01%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "229.59"
$ws.Range("E22").Value = "  +2.02%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "2.19"
$ws.Range("E23").Value = "  +13.11%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "9.32"
$ws.Range("E24").Value = "  -2.96%  "

$ws.Range("E25").Value = "  +0.09%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "10.70"
$ws.Range("E26").Value = "  +1.38%  "

$ws.Range("E27").Value = "  +3.35%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "38.81"
$ws.Range("E28").Value = "  +24.74%  "

$ws.Range("B29").Value = "Toncoin"
$ws.Range("C29").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.24"
$ws.Range("E29").Value = "  +3.57%  "

$ws.Range("B30").Value = "PancakeSwap"
$ws.Range("C30").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.22"
$ws.Range("E30").Value = "  +3.23%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "170.25"
$ws.Range("E31").Value = "  +0.00%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "20.17"
$ws.Range("E32").Value = "  +2.90%  "

$ws.Range("E33").Value = "  +6.03%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.24"
$ws.Range("E34").Value = "  +3.26%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.114"
$ws.Range("E35").Value = "  +11.33%  "

$ws.Range("E36").Value = "  +1.16%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.47"
$ws.Range("E37").Value = "  +6.49%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0326"
$ws.Range("E38").Value = "  +12.16%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "12.35"
$ws.Range("E39").Value = "  +5.41%  "

$ws.Range("E40").Value = "  +3.18%  "

$ws.Range("E41").Value = "  +9.74%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.36"
$ws.Range("E42").Value = "  +2.70%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "59.59"
$ws.Range("E43").Value = "  +3.34%  "

$ws.Range("E44").Value = "  +29.47%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "8.65"
$ws.Range("E45").Value = "  +5.60%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "103.08"
$ws.Range("E46").Value = "  +6.08%  "

$ws.Range("E47").Value = "  +3.11%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.43"
$ws.Range("E48").Value = "  +13.71%  "

$ws.Range("E49").Value = "  +3.70%  "

$ws.Range("E50").Value = "  +2.89%  "

$ws.Range("E51").Value = "  +2.37%  "
